$d = $word.ActiveDocument

# PAR1 boundary 0: "1) Princípio de Stev..." -> "2) Empuxo e Princípi..."
$r = $d.Content
$found = $r.Find.Execute("vin e Pascal2) Empuxo e Prin", $false, $false, $false, $false, $false, $true, 1, $false, "vin e Pascal^l2) Empuxo e Prin", 2)
if (-not $found) { throw "Not found: PAR1 boundary 0" }

# PAR1 boundary 1: "2) Empuxo e Princípi..." -> "3) Tensão superficia..."
$r = $d.Content
$found = $r.Find.Execute("e Arquimedes3) Tensão superf", $false, $false, $false, $false, $false, $true, 1, $false, "e Arquimedes^l3) Tensão superf", 2)
if (-not $found) { throw "Not found: PAR1 boundary 1" }

# PAR1 boundary 2: "3) Tensão superficia..." -> "4) Queda em um meio ..."
$r = $d.Content
$found = $r.Find.Execute(" superficial4) Queda em um m", $false, $false, $false, $false, $false, $true, 1, $false, " superficial^l4) Queda em um m", 2)
if (-not $found) { throw "Not found: PAR1 boundary 2" }

# PAR1 boundary 3: "4) Queda em um meio ..." -> "5) Sistema massa-mol..."
$r = $d.Content
$found = $r.Find.Execute("meio viscoso5) Sistema massa", $false, $false, $false, $false, $false, $true, 1, $false, "meio viscoso^l5) Sistema massa", 2)
if (-not $found) { throw "Not found: PAR1 boundary 3" }

# PAR1 boundary 4: "5) Sistema massa-mol..." -> "6) Ondas mecânicas..."
$r = $d.Content
$found = $r.Find.Execute("a massa-mola6) Ondas mecânic", $false, $false, $false, $false, $false, $true, 1, $false, "a massa-mola^l6) Ondas mecânic", 2)
if (-not $found) { throw "Not found: PAR1 boundary 4" }

# PAR1 boundary 5: "6) Ondas mecânicas..." -> "7) Calor, temperatur..."
$r = $d.Content
$found = $r.Find.Execute("as mecânicas7) Calor, temper", $false, $false, $false, $false, $false, $true, 1, $false, "as mecânicas^l7) Calor, temper", 2)
if (-not $found) { throw "Not found: PAR1 boundary 5" }

# PAR1 boundary 6: "7) Calor, temperatur..." -> "8) Dilatação linear..."
$r = $d.Content
$found = $r.Find.Execute("enar energia8) Dilatação lin", $false, $false, $false, $false, $false, $true, 1, $false, "enar energia^l8) Dilatação lin", 2)
if (-not $found) { throw "Not found: PAR1 boundary 6" }

# PAR1 boundary 7: "8) Dilatação linear..." -> "9) Os meios de propa..."
$r = $d.Content
$found = $r.Find.Execute("tação linear9) Os meios de p", $false, $false, $false, $false, $false, $true, 1, $false, "tação linear^l9) Os meios de p", 2)
if (-not $found) { throw "Not found: PAR1 boundary 7" }

# PAR1 boundary 8: "9) Os meios de propa..." -> "10) Calor específico..."
$r = $d.Content
$found = $r.Find.Execute("ção de calor10) Calor especí", $false, $false, $false, $false, $false, $true, 1, $false, "ção de calor^l10) Calor especí", 2)
if (-not $found) { throw "Not found: PAR1 boundary 8" }

# PAR1 boundary 9: "10) Calor específico..." -> "11) A lei de Boyle-M..."
$r = $d.Content
$found = $r.Find.Execute("alor latente11) A lei de Boy", $false, $false, $false, $false, $false, $true, 1, $false, "alor latente^l11) A lei de Boy", 2)
if (-not $found) { throw "Not found: PAR1 boundary 9" }

# PAR2 boundary 0: "1) stevin’s and Pasc..." -> "2) Thrust and Archim..."
$r = $d.Content
$found = $r.Find.Execute("’s Principle2) Thrust and Ar", $false, $false, $false, $false, $false, $true, 1, $false, "’s Principle^l2) Thrust and Ar", 2)
if (-not $found) { throw "Not found: PAR2 boundary 0" }

# PAR2 boundary 1: "2) Thrust and Archim..." -> "3) Surface tension..."
$r = $d.Content
$found = $r.Find.Execute("s’ Principle3) Surface tensi", $false, $false, $false, $false, $false, $true, 1, $false, "s’ Principle^l3) Surface tensi", 2)
if (-not $found) { throw "Not found: PAR2 boundary 1" }

# PAR2 boundary 2: "3) Surface tension..." -> "4) The fall in a vis..."
$r = $d.Content
$found = $r.Find.Execute("face tension4) The fall in a", $false, $false, $false, $false, $false, $true, 1, $false, "face tension^l4) The fall in a", 2)
if (-not $found) { throw "Not found: PAR2 boundary 2" }

# PAR2 boundary 3: "4) The fall in a vis..." -> "5) Mass-spring syste..."
$r = $d.Content
$found = $r.Find.Execute("iscous fluid5) Mass-spring s", $false, $false, $false, $false, $false, $true, 1, $false, "iscous fluid^l5) Mass-spring s", 2)
if (-not $found) { throw "Not found: PAR2 boundary 3" }

# PAR2 boundary 4: "5) Mass-spring syste..." -> "6) Mechanical waves..."
$r = $d.Content
$found = $r.Find.Execute("pring system6) Mechanical wa", $false, $false, $false, $false, $false, $true, 1, $false, "pring system^l6) Mechanical wa", 2)
if (-not $found) { throw "Not found: PAR2 boundary 4" }

# PAR2 boundary 5: "6) Mechanical waves..." -> "7) Heat, temperature..."
$r = $d.Content
$found = $r.Find.Execute("anical waves7) Heat, tempera", $false, $false, $false, $false, $false, $true, 1, $false, "anical waves^l7) Heat, tempera", 2)
if (-not $found) { throw "Not found: PAR2 boundary 5" }

# PAR2 boundary 6: "7) Heat, temperature..." -> "8) Linear thermal ex..."
$r = $d.Content
$found = $r.Find.Execute("store energy8) Linear therma", $false, $false, $false, $false, $false, $true, 1, $false, "store energy^l8) Linear therma", 2)
if (-not $found) { throw "Not found: PAR2 boundary 6" }

# PAR2 boundary 7: "8) Linear thermal ex..." -> "9) The fundamental m..."
$r = $d.Content
$found = $r.Find.Execute("al expansion9) The fundament", $false, $false, $false, $false, $false, $true, 1, $false, "al expansion^l9) The fundament", 2)
if (-not $found) { throw "Not found: PAR2 boundary 7" }

# PAR2 boundary 8: "9) The fundamental m..." -> "10) Specific and lat..."
$r = $d.Content
$found = $r.Find.Execute("eat transfer10) Specific and", $false, $false, $false, $false, $false, $true, 1, $false, "eat transfer^l10) Specific and", 2)
if (-not $found) { throw "Not found: PAR2 boundary 8" }

# PAR2 boundary 9: "10) Specific and lat..." -> "11) The Boyle-Mariot..."
$r = $d.Content
$found = $r.Find.Execute(" latent heat11) The Boyle-Ma", $false, $false, $false, $false, $false, $true, 1, $false, " latent heat^l11) The Boyle-Ma", 2)
if (-not $found) { throw "Not found: PAR2 boundary 9" }

# PAR3 boundary 0: "1. Apostilas do Labo..." -> "2. VUOLO, J.H. Funda..."
$r = $d.Content
$found = $r.Find.Execute("do IFSC/USP.2. VUOLO, J.H. F", $false, $false, $false, $false, $false, $true, 1, $false, "do IFSC/USP.^l2. VUOLO, J.H. F", 2)
if (-not $found) { throw "Not found: PAR3 boundary 0" }

# PAR3 boundary 1: "2. VUOLO, J.H. Funda..." -> "3. NUSSENZVEIG, H.M...."
$r = $d.Content
$found = $r.Find.Execute("cher (1996).3. NUSSENZVEIG, ", $false, $false, $false, $false, $false, $true, 1, $false, "cher (1996).^l3. NUSSENZVEIG, ", 2)
if (-not $found) { throw "Not found: PAR3 boundary 1" }

# PAR3 boundary 2: "3. NUSSENZVEIG, H.M...." -> "4. RESNICK, R.; HALL..."
$r = $d.Content
$found = $r.Find.Execute("cher (2008).4. RESNICK, R.; ", $false, $false, $false, $false, $false, $true, 1, $false, "cher (2008).^l4. RESNICK, R.; ", 2)
if (-not $found) { throw "Not found: PAR3 boundary 2" }

# PAR3 boundary 3: "4. RESNICK, R.; HALL..." -> "5. TIPLER, P.; MOSCA..."
$r = $d.Content
$found = $r.Find.Execute(" LTC (2008).5. TIPLER, P.; M", $false, $false, $false, $false, $false, $true, 1, $false, " LTC (2008).^l5. TIPLER, P.; M", 2)
if (-not $found) { throw "Not found: PAR3 boundary 3" }

# PAR3 boundary 4: "5. TIPLER, P.; MOSCA..." -> "6. SEARS, F. W.; ZEM..."
$r = $d.Content
$found = $r.Find.Execute(" LTC (2008).6. SEARS, F. W.;", $false, $false, $false, $false, $false, $true, 1, $false, " LTC (2008).^l6. SEARS, F. W.;", 2)
if (-not $found) { throw "Not found: PAR3 boundary 4" }

# PAR3 boundary 5: "6. SEARS, F. W.; ZEM..." -> "    Pearson Addison ..."
$r = $d.Content
$found = $r.Find.Execute("II, Vol. 2,     Pearson Addi", $false, $false, $false, $false, $false, $true, 1, $false, "II, Vol. 2, ^l    Pearson Addi", 2)
if (-not $found) { throw "Not found: PAR3 boundary 5" }

# PAR3 boundary 6: "    Pearson Addison ..." -> "7. JEWETT Jr, John W..."
$r = $d.Content
$found = $r.Find.Execute("sley (2009).7. JEWETT Jr, Jo", $false, $false, $false, $false, $false, $true, 1, $false, "sley (2009).^l7. JEWETT Jr, Jo", 2)
if (-not $found) { throw "Not found: PAR3 boundary 6" }

Write-Output "All replacements applied successfully."
